$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing descriptions (text content changes, same row/cell) ---

# 隐身 (Stealth) description: add "持续时间结束或" before "自身攻击后退出隐身状态"
$ws.Range("C3").Value = "隐身状态下，敌方无法选中你作为目标，持续时间结束或自身攻击后退出隐身状态"

# 流血 (Bleed) description: now triggers on taking damage (extra damage equal to stacks),
# rather than at the start of the turn. The old wording is reused below for the new "腐蚀" row.
$ws.Range("C4").Value = "拥有流血状态的角色受到伤害时，额外受到流血层数的伤害"

# --- Append new skill/status rows ---

# Row 6: 念力 (Psychic Force)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "念力"
$ws.Range("C6").Value = "拥有念力的角色在造成伤害时，额外造成念力层数的伤害"

# Row 7: 腐蚀 (Corrosion) - reuses the original 流血 wording
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "腐蚀"
$ws.Range("C7").Value = "拥有流血状态的角色回合开始时受到一次流血层数的伤害"

# Row 8: 额外回合 (Extra Turn)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "额外回合"
$ws.Range("C8").Value = "不消耗buff持续时间的特殊回合"

# --- Update the sheet's active selection to match the new end of data ---
$null = $ws.Range("C10").Select()
